# Applies the cell-value corrections from the scheduled-runner price refresh.
# Each worksheet ("ALC","ARM","BSM","CRP","CUL","GSM","LTW","WVR") stores a
# leve-crafting profit table; columns H-N are computed market-price figures
# that the runner recomputed. We just overwrite the stored values directly,
# matching how the source data was produced (no live formulas in this sheet).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 5666.5576
$ws.Range("J17").Value = 6029.8423
$ws.Range("L17").Value = 18089.5269
$ws.Range("N17").Value = -18425.5269
$ws.Range("H18").Value = 2718.7896
$ws.Range("I18").Value = 2670.3333
$ws.Range("K18").Value = 2670.3333
$ws.Range("M18").Value = -2386.3333
$ws.Range("H28").Value = 350.8125
$ws.Range("I28").Value = 239.61539
$ws.Range("K28").Value = 239.61539
$ws.Range("M28").Value = 245.38461
$ws.Range("H62").Value = 4269.3335
$ws.Range("I62").Value = 4152.5
$ws.Range("K62").Value = 4152.5
$ws.Range("M62").Value = -3528.5
$ws.Range("H65").Value = 4269.3335
$ws.Range("I65").Value = 4152.5
$ws.Range("K65").Value = 20762.5
$ws.Range("M65").Value = -17642.5
$ws.Range("H98").Value = 2821.8696
$ws.Range("I98").Value = 2138.2856
$ws.Range("K98").Value = 2138.2856
$ws.Range("M98").Value = -640.2856000000002
$ws.Range("H122").Value = 2821.8696
$ws.Range("I122").Value = 2138.2856
$ws.Range("K122").Value = 6414.8568
$ws.Range("M122").Value = -3964.8568
$ws.Range("H132").Value = 4524.92
$ws.Range("I132").Value = 4874.9473
$ws.Range("K132").Value = 14624.8419
$ws.Range("M132").Value = -12094.8419
$ws.Range("H138").Value = 2836.8518
$ws.Range("I138").Value = 1320.64
$ws.Range("K138").Value = 3961.92
$ws.Range("M138").Value = 1178.08

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2815.2144
$ws.Range("I61").Value = 2492.5454
$ws.Range("J61").Value = 3998.3333
$ws.Range("K61").Value = 2492.5454
$ws.Range("L61").Value = 3998.3333
$ws.Range("M61").Value = -2280.5454
$ws.Range("N61").Value = -4422.3333
$ws.Range("H74").Value = 62987.062
$ws.Range("I74").Value = 7029.0356
$ws.Range("J74").Value = 454693.25
$ws.Range("K74").Value = 7029.0356
$ws.Range("L74").Value = 454693.25
$ws.Range("M74").Value = -6155.0356
$ws.Range("N74").Value = -456441.25
$ws.Range("H77").Value = 62987.062
$ws.Range("I77").Value = 7029.0356
$ws.Range("J77").Value = 454693.25
$ws.Range("K77").Value = 35145.178
$ws.Range("L77").Value = 2273466.25
$ws.Range("M77").Value = -30777.178
$ws.Range("N77").Value = -2282202.25
$ws.Range("H122").Value = 2201.8057
$ws.Range("I122").Value = 1966.3928
$ws.Range("J122").Value = 3025.75
$ws.Range("K122").Value = 5899.178400000001
$ws.Range("L122").Value = 9077.25
$ws.Range("M122").Value = -3449.178400000001
$ws.Range("N122").Value = -13977.25
$ws.Range("H132").Value = 2583.375
$ws.Range("I132").Value = 1659.2
$ws.Range("K132").Value = 4977.6
$ws.Range("M132").Value = -2447.6
$ws.Range("H136").Value = 2815.2144
$ws.Range("I136").Value = 2492.5454
$ws.Range("J136").Value = 3998.3333
$ws.Range("K136").Value = 7477.6362
$ws.Range("L136").Value = 11994.9999
$ws.Range("M136").Value = -4927.6362
$ws.Range("N136").Value = -17094.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H103").Value = 20000
$ws.Range("J103").Value = 20000
$ws.Range("L103").Value = 20000
$ws.Range("N103").Value = -22344
$ws.Range("H105").Value = 1643.7333
$ws.Range("I105").Value = 1743.4615
$ws.Range("K105").Value = 1743.4615
$ws.Range("M105").Value = 3.538500000000113
$ws.Range("H134").Value = 3289.5642
$ws.Range("I134").Value = 1587.3
$ws.Range("K134").Value = 4761.9
$ws.Range("M134").Value = -2226.9

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 16355.187
$ws.Range("I31").Value = 1594.6364
$ws.Range("J31").Value = 18892.156
$ws.Range("K31").Value = 1594.6364
$ws.Range("L31").Value = 18892.156
$ws.Range("M31").Value = -1299.6364
$ws.Range("N31").Value = -19482.156
$ws.Range("H34").Value = 16355.187
$ws.Range("I34").Value = 1594.6364
$ws.Range("J34").Value = 18892.156
$ws.Range("K34").Value = 1594.6364
$ws.Range("L34").Value = 18892.156
$ws.Range("M34").Value = -1392.6364
$ws.Range("N34").Value = -19296.156
$ws.Range("H58").Value = 2574.8462
$ws.Range("I58").Value = 1716.8889
$ws.Range("K58").Value = 1716.8889
$ws.Range("M58").Value = -1513.8889
$ws.Range("H62").Value = 2325
$ws.Range("I62").Value = 1766.6666
$ws.Range("K62").Value = 1766.6666
$ws.Range("M62").Value = -1142.6666
$ws.Range("H65").Value = 2325
$ws.Range("I65").Value = 1766.6666
$ws.Range("K65").Value = 8833.333000000001
$ws.Range("M65").Value = -5713.333000000001
$ws.Range("H94").Value = 1193.6154
$ws.Range("I94").Value = 928
$ws.Range("J94").Value = 1273.3
$ws.Range("K94").Value = 928
$ws.Range("L94").Value = 1273.3
$ws.Range("M94").Value = -477
$ws.Range("N94").Value = -2175.3
$ws.Range("H105").Value = 2399.2
$ws.Range("I105").Value = 2261.5
$ws.Range("J105").Value = 2950
$ws.Range("K105").Value = 2261.5
$ws.Range("L105").Value = 2950
$ws.Range("M105").Value = -514.5
$ws.Range("N105").Value = -6444
$ws.Range("H132").Value = 32152.4
$ws.Range("I132").Value = 2486.8462
$ws.Range("J132").Value = 224978.5
$ws.Range("K132").Value = 7460.5386
$ws.Range("L132").Value = 674935.5
$ws.Range("M132").Value = -4930.5386
$ws.Range("N132").Value = -679995.5
$ws.Range("H136").Value = 2574.8462
$ws.Range("I136").Value = 1716.8889
$ws.Range("K136").Value = 5150.6667
$ws.Range("M136").Value = -2600.6667

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 166.85715
$ws.Range("J23").Value = 163
$ws.Range("L23").Value = 489
$ws.Range("N23").Value = -959
$ws.Range("H39").Value = 653.5714
$ws.Range("J39").Value = 815
$ws.Range("L39").Value = 2445
$ws.Range("N39").Value = -3033
$ws.Range("H47").Value = 543.4286
$ws.Range("I47").Value = 600.6667
$ws.Range("J47").Value = 200
$ws.Range("K47").Value = 1802.0001
$ws.Range("L47").Value = 600
$ws.Range("M47").Value = -1371.0001
$ws.Range("N47").Value = -1462
$ws.Range("H56").Value = 16672133
$ws.Range("I56").Value = 16672133
$ws.Range("K56").Value = 16672133
$ws.Range("M56").Value = -16671603

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 1085
$ws.Range("I31").Value = 1085
$ws.Range("K31").Value = 1085
$ws.Range("M31").Value = -793
$ws.Range("H37").Value = 1085
$ws.Range("I37").Value = 1085
$ws.Range("K37").Value = 1085
$ws.Range("M37").Value = -808
$ws.Range("H46").Value = 9849.143
$ws.Range("J46").Value = 19997.5
$ws.Range("L46").Value = 19997.5
$ws.Range("N46").Value = -20309.5
$ws.Range("H113").Value = 1516.3334
$ws.Range("J113").Value = 599
$ws.Range("L113").Value = 599
$ws.Range("N113").Value = -4939
$ws.Range("H126").Value = 3182.625
$ws.Range("I126").Value = 2944.45
$ws.Range("K126").Value = 8833.349999999999
$ws.Range("M126").Value = -6363.349999999999
$ws.Range("H132").Value = 3571.1904
$ws.Range("I132").Value = 2830.1667
$ws.Range("J132").Value = 4559.222
$ws.Range("K132").Value = 8490.500100000001
$ws.Range("L132").Value = 13677.666
$ws.Range("M132").Value = -5960.500100000001
$ws.Range("N132").Value = -18737.666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 15800.375
$ws.Range("I32").Value = 1599.75
$ws.Range("J32").Value = 30001
$ws.Range("K32").Value = 1599.75
$ws.Range("L32").Value = 30001
$ws.Range("M32").Value = -1282.75
$ws.Range("N32").Value = -30635
$ws.Range("H46").Value = 5121.3335
$ws.Range("I46").Value = 4293.5
$ws.Range("J46").Value = 5783.6
$ws.Range("K46").Value = 4293.5
$ws.Range("L46").Value = 5783.6
$ws.Range("M46").Value = -4105.5
$ws.Range("N46").Value = -6159.6
$ws.Range("H61").Value = 2804.5
$ws.Range("I61").Value = 2804.5
$ws.Range("K61").Value = 2804.5
$ws.Range("M61").Value = -2602.5
$ws.Range("H113").Value = 2804.5
$ws.Range("I113").Value = 2804.5
$ws.Range("K113").Value = 2804.5
$ws.Range("M113").Value = -634.5
$ws.Range("H136").Value = 29064.205
$ws.Range("I136").Value = 40950.42
$ws.Range("J136").Value = 5291.769
$ws.Range("K136").Value = 122851.26
$ws.Range("L136").Value = 15875.307
$ws.Range("M136").Value = -120301.26
$ws.Range("N136").Value = -20975.307
$ws.Range("H139").Value = 75000
$ws.Range("J139").Value = 75000
$ws.Range("L139").Value = 75000
$ws.Range("N139").Value = -85280

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 4576.25
$ws.Range("I17").Value = 2000
$ws.Range("J17").Value = 5435
$ws.Range("K17").Value = 2000
$ws.Range("L17").Value = 5435
$ws.Range("M17").Value = -1828
$ws.Range("N17").Value = -5779
$ws.Range("H113").Value = 1285.2858
$ws.Range("I113").Value = 999.5
$ws.Range("J113").Value = 1399.6
$ws.Range("K113").Value = 2998.5
$ws.Range("L113").Value = 4198.799999999999
$ws.Range("M113").Value = -828.5
$ws.Range("N113").Value = -8538.799999999999
$ws.Range("H122").Value = 1601.8148
$ws.Range("I122").Value = 1054.5
$ws.Range("K122").Value = 3163.5
$ws.Range("M122").Value = -713.5
$ws.Range("H132").Value = 15437.286
$ws.Range("I132").Value = 2787.8667
$ws.Range("K132").Value = 8363.6001
$ws.Range("M132").Value = -5833.6001
$ws.Range("H136").Value = 3668.1538
$ws.Range("I136").Value = 3198.7273
$ws.Range("J136").Value = 6250
$ws.Range("K136").Value = 9596.1819
$ws.Range("L136").Value = 18750
$ws.Range("M136").Value = -7046.1819
$ws.Range("N136").Value = -23850
